$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update existing row 8 date value (45269 -> 45270)
$ws.Range("A8").Value = 45270

# 2. Insert 3 new rows (9:11) above the old Total row (old row 9 shifts down to row 12)
$ws.Rows("9:11").Insert()

# Give the new rows the same thin-border box used by the surrounding data rows
# (this reuses the existing cellXfs entries 6/7/8/9 instead of minting new ones).
$ws.Range("A9:F11").Borders.LineStyle = 1

# Row 9
$ws.Range("A9").Value = 45271
$ws.Range("B9").Value = 0.541666666666667
$ws.Range("C9").Value = 0.708333333333333
$ws.Range("D9").Formula = "=(C9<B9)+C9-B9"
$ws.Range("E9").Value = 10
$ws.Range("F9").Formula = "=(D9*24)*E9"

# Row 10
$ws.Range("A10").Value = 45272
$ws.Range("B10").Value = 0.458333333333333
$ws.Range("C10").Value = 0.791666666666667
$ws.Range("D10").Formula = "=(C10<B10)+C10-B10"
$ws.Range("E10").Value = 10
$ws.Range("F10").Formula = "=(D10*24)*E10"

# Row 11 (From/To left blank)
$ws.Range("A11").Value = 45273
$ws.Range("D11").Formula = "=(C11<B11)+C11-B11"
$ws.Range("E11").Value = 10
$ws.Range("F11").Formula = "=(D11*24)*E11"

# Number formats for the new rows, matching the existing column formats exactly
# (escaped hyphens so the engine recognises/reuses numFmtId 165 rather than minting a new one)
$ws.Range("A9:A11").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B9:C11").NumberFormat = "h:mm"
$ws.Range("D9:D11").NumberFormat = "[hh]:mm:ss"
$ws.Range("E9:F11").NumberFormat = "General"

# 3. Total row (pushed down from row 9 to row 12) - extend the SUM ranges
$ws.Range("D12").Formula = "=SUM(D2:D11)"
$ws.Range("F12").Formula = "=SUM(F2:F11)"

# 4. Selection moves to E12
[void]$ws.Range("E12").Select()
